# AERSP 424 Project Report - apply commit "changes and conclusions to the report"
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
    return $found
}

# 1) "The draft constructor ..." paragraph rewrite
Replace-Text `
    "The draft constructor takes in the league members, position limits, round limit, and draft type. In its body, the constructor fills the team player map with the league members, and then chooses which CSV file to open based on the user input. It then resizes the positional limit containers depending on the number of league members. " `
    "The draft constructor receives the inputs of league members, position limits, round limit, and draft type. In its body, the constructor fills the team player map with the league members, and then chooses which CSV file to open based on the user input. The team player map is used throughout the code to keep track of the team and their corresponding players and positions. It then resizes the positional limit containers depending on the number of league members. "

# 2) ". The first instance displays ..." -> add "is considered the default application and"
Replace-Text `
    ". The first instance displays the top ten players available to the user. The second instance allows the user to input the number of players they would like to see, and then outputs that number of players. " `
    ". The first instance is considered the default application and displays the top ten players available to the user. The second instance allows the user to input the number of players they would like to see, and then outputs that number of players. "

# 3) "Next up is the getUserPick function ..." text is unchanged in the commit (only a run
#    split / lastRenderedPageBreak marker move happens in Word's internal layout) -> no text edit needed

# 4) "This function simulates the computer's pick ..." paragraph rewrite
Replace-Text `
    "This function simulates the computer’s pick for the other members of the draft, as well as calling the getUserPick function to have the user make their own pick. The computer makes its pick through using probability rules which I will describe later. If the player chosen by the computer plays a position that has already met its limit, the code will output “Moving to next player” and draft a player that will not exceed the limit. If the CSV file runs out of players needed to fill a position, the computer will forfeit that pick and say “No more of the position desired available, voiding pick”. Another cool feature of this function is the ability to have a snake draft. A snake draft is where the draft order switches at the start of every other round. For example, if John picked last in the first round, he would pick first in the second round. This is done by checking if the round number is divisible evenly by two at the start of every round. " `
    "This function simulates the computer’s pick for the other members of the draft, as well as calling the getUserPick function to have the user make their own pick. The computer makes its pick through using probability rules which will be described later. If the player chosen by the computer plays a position that has already met its limit, the code will output “Moving to next player” and draft a player that will not exceed the limit. If the CSV file runs out of players needed to fill a position, the computer will forfeit that pick and say “No more of the position desired available, voiding pick”. Another cool feature of this function is the ability to have a snake draft. A snake draft is where the draft order switches at the start of every other round. For example, if a team picked last in the first round, they would pick first in the second round. This is done by checking if the round number is divisible evenly by two at the start of every round. "

# 5) "The computer's pick probability was done ..." paragraph rewrite
Replace-Text `
    "The computer’s pick probability was done through a function called pickRandomizer. This function aligns the top four players available with the probabilities 40%, 30%, 20%, and 10% respectively. The player is then picked based on these percentages. This was done to represent some of the chaos that may occur in a real-life draft where the top player available will not always be chosen. Next the picks on each team are displayed with the outputTeamPicks function. This function iterates through every player and outputs their picks in order of QB, RB, WR, and TE. This function also keeps track of the average salary of a team which may be important to the user. Finally, the last function in this class is the operate function. The operate function gives the user the choice to automatically draft their picks for them if they choose. It also handles the case where the draft file cannot be opened for any reason. Lastly, this function calls all other important functions in the class and allows the draft to run smoothly. " `
    "The computer’s pick was done with a probability method to add randomness to the draft to simulate a real-life draft scenario. The logic for the probability implementation was completed with a function called pickRandomizer. This function aligns the top four players available with the probabilities 40%, 30%, 20%, and 10% respectively. The player is then picked based on these percentages. This was done to represent some of the chaos that may occur in a real-life draft where the top player available will not always be chosen. Next the picks on each team are displayed with the outputTeamPicks function. This function iterates through every player and outputs their picks in order of QB, RB, WR, and TE. This function also keeps track of the estimated average salary of a team which may be important or interesting to the user. Finally, the last function in this class is the operate function. The operate function gives the user the choice to automatically draft their picks for them if they choose. It also handles the case where the draft file cannot be opened for any reason. Lastly, this function calls all other important functions in the class and allows the draft to run smoothly. "

# 6) "Based on what the user inputs, the function instantiates ..." paragraph rewrite
#    (also removes the lastRenderedPageBreak that used to sit mid-paragraph here)
Replace-Text `
    "Based on what the user inputs, the function instantiates an object from each of the ranking, CustomLeague, and draft classes. Additionally, in the CustomLeague object, a thread is created to run the functionality of the CustomLeague. After this, the user is asked if they want to complete another draft and if so which type. The process I just described is done for all three draft types (standard, ppr, and half ppr). This means the ranking, CustomLeague, and draft classes all have three objects instantiated from them. Finally, the number of drafts completed is shown to the user, and the code has finished. " `
    "Based on what the user inputs, the function instantiates an object from each of the “Ranking,” “CustomLeague,” and “Draft” classes. Additionally, in the CustomLeague object, a thread is created to run the functionality of the CustomLeague. After this, the user is asked if they want to complete another draft and if so which type. The process described is done for all three draft types (standard, ppr, and half ppr). This means the “Ranking,” “CustomLeague,” and “Draft” classes all have three objects instantiated from them. Finally, the number of drafts completed is shown to the user, and the code has finished. "

# 7) "Add some conclusions about the code" -> full conclusions paragraph
Replace-Text `
    "Add some conclusions about the code" `
    "The code in this project provides a unique drafting experience for a user. While there are platforms that exist to provide mock drafts, few have quite the capability that this code provides. Some unique features include position limits and full flexibility of team size, number of league members, and number of drafts completed. In addition, this code can work with any CSV file with the same format as the files used in this project. Therefore, any user can create a CSV file of rankings and draft based on those rankings. Few, if any, platforms provide that capability. In summary, this project provides a user with a great mock drafting experience to prepare them for the upcoming fantasy football league year.  "
